$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: add status "en proceso" in column C
$ws.Range("C7").Value = "en proceso"

# Row 8: change the raw number 100 into a percentage-formatted 1 (100%)
$ws.Range("C8").Value = 1
$ws.Range("C8").NumberFormat = "0%"

# Row 11: add an empty, underlined placeholder cell in column D
$ws.Range("D11").Font.Underline = $true

# Row 12: add "Cancelada" status and explanation text
# (set D12 before C12 so new shared-string entries are appended in the
#  same order as the original edit: "No, esta tarea..." then "Cancelada")
$ws.Range("D12").Value = "No, esta tarea la suprimimos porque es la opción de que un cliente modifique el pago de sus cuotas"
$ws.Range("C12").Value = "Cancelada"

# Row 13: add an empty, underlined placeholder cell in column C
$ws.Range("C13").Font.Underline = $true

# Row 17: replace the text status with a percentage-formatted 1 (100%)
$ws.Range("C17").Value = 1
$ws.Range("C17").NumberFormat = "0%"

# Row 18: add an empty, underlined placeholder cell in column C
$ws.Range("C18").Font.Underline = $true

# Row 26 (new row): new task and its owner
$ws.Range("A26").Value = "Agregar patron fechas en todos los campos de fecha"
$ws.Range("B26").Value = "Agustina"

# Update the active selection to match the edited area
$ws.Range("B25").Select()
